# Generate Report for Handoff
# Adds two new handed-off files (21350e15-...-md and 4426d585-...-md) as
# new rows to the Overview sheet and to the per-locale (zh-cn / de-de)
# detail sheets, expanding each sheet's table accordingly.

$wb = $excel.ActiveWorkbook

$file1Name  = "21350e15-1436-418f-bf98-2c2f74d1a38c.md"
$file1Path  = "e2e\21350e15-1436-418f-bf98-2c2f74d1a38c.md"
$file2Name  = "4426d585-d49b-440d-8e9c-305b36e872d1.md"
$file2Path  = "e2e\4426d585-d49b-440d-8e9c-305b36e872d1.md"

$overviewDate = "2016-08-21 18:45:31"

$zhXlf1 = "21350e15-1436-418f-bf98-2c2f74d1a38c.fc4a7cf1a146852f73246e09337b6ca002106577.zh-cn.xlf"
$zhXlf2 = "4426d585-d49b-440d-8e9c-305b36e872d1.9b2c0e3b6c45ffdeef358085f389b891e2884f0d.zh-cn.xlf"
$zhDate = "2016-08-21 18:45:27"

$deXlf1 = "21350e15-1436-418f-bf98-2c2f74d1a38c.fc4a7cf1a146852f73246e09337b6ca002106577.de-de.xlf"
$deXlf2 = "4426d585-d49b-440d-8e9c-305b36e872d1.9b2c0e3b6c45ffdeef358085f389b891e2884f0d.de-de.xlf"
$deDate = "2016-08-21 18:45:31"

$srcHref1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f711e842bf88c9fbb379ec8a3f81c37e88d1105b/e2e/21350e15-1436-418f-bf98-2c2f74d1a38c.md"
$srcHref2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5c8c521e7f81ce288c16fab0e366420f33c57356/e2e/4426d585-d49b-440d-8e9c-305b36e872d1.md"

# ---------------------------------------------------------------------
# Sheet "Overview" (sheet1): File Name | Path And Name | Extension |
#   Publish URL | zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

$loOverview.ListRows.Add() | Out-Null
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $file1Name
$wsOverview.Range("B4").Value = $file1Path
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = ""
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = $overviewDate
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Range("A5").Value = $file2Name
$wsOverview.Range("B5").Value = $file2Path
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = ""
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = $overviewDate
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $srcHref1, "", "", $file1Path) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), $srcHref2, "", "", $file2Path) | Out-Null

$wsOverview.Range("B4").Style = "HyperLink"
$wsOverview.Range("B5").Style = "HyperLink"

# ---------------------------------------------------------------------
# Sheet "zh-cn" (sheet2): Source File Name | File Extension | Status |
#   Source Path | Priority | Content Duplicate | Latest Handoff File |
#   Latest Handoff Datetime | Latest Target File | Latest Handback File |
#   Latest Handback DateTime | Reference Tokens | To be localized |
#   Dependency From | Has metadata | Error Detail
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)

$loZh.ListRows.Add() | Out-Null
$loZh.ListRows.Add() | Out-Null

function Fill-DetailRow($ws, [int]$r, [string]$fileName, [string]$xlfName, [string]$xlfDate) {
    $ws.Range("A$r").Value = $fileName
    $ws.Range("B$r").Value = ".md"
    $ws.Range("C$r").Value = "Ready for handoff"
    $ws.Range("D$r").Value = "e2e"
    $ws.Range("E$r").Value = "ht"
    # Leading apostrophe forces text interpretation so "True"/"False" are
    # stored as literal text (shared string), matching the source workbook,
    # instead of being auto-coerced to native boolean cells.
    $ws.Range("F$r").Value = "'False"
    $ws.Range("G$r").Value = $xlfName
    $ws.Range("H$r").Value = $xlfDate
    $ws.Range("H$r").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("I$r").Value = ""
    $ws.Range("J$r").Value = ""
    $ws.Range("K$r").Value = "0001-01-01 00:00:00"
    $ws.Range("K$r").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("L$r").Value = ""
    $ws.Range("M$r").Value = "'True"
    $ws.Range("N$r").Value = ""
    $ws.Range("O$r").Value = "'False"
    $ws.Range("P$r").Value = ""
}

Fill-DetailRow $wsZh 4 $file1Name $zhXlf1 $zhDate
Fill-DetailRow $wsZh 5 $file2Name $zhXlf2 $zhDate

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $srcHref1, "", "", $file1Name) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), $srcHref2, "", "", $file2Name) | Out-Null

$wsZh.Range("A4").Style = "HyperLink"
$wsZh.Range("A5").Style = "HyperLink"

# ---------------------------------------------------------------------
# Sheet "de-de" (sheet3): same columns as "zh-cn"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)

$loDe.ListRows.Add() | Out-Null
$loDe.ListRows.Add() | Out-Null

Fill-DetailRow $wsDe 4 $file1Name $deXlf1 $deDate
Fill-DetailRow $wsDe 5 $file2Name $deXlf2 $deDate

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $srcHref1, "", "", $file1Name) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), $srcHref2, "", "", $file2Name) | Out-Null

$wsDe.Range("A4").Style = "HyperLink"
$wsDe.Range("A5").Style = "HyperLink"
